$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row pairs (columns F:V; A:E stay fixed to row position) ---
$rowA = $ws.Range("F51:V51").Value2
$rowB = $ws.Range("F52:V52").Value2
$ws.Range("F51:V51").Value2 = $rowB
$ws.Range("F52:V52").Value2 = $rowA

$rowA = $ws.Range("F53:V53").Value2
$rowB = $ws.Range("F54:V54").Value2
$ws.Range("F53:V53").Value2 = $rowB
$ws.Range("F54:V54").Value2 = $rowA

$rowA = $ws.Range("F72:V72").Value2
$rowB = $ws.Range("F73:V73").Value2
$ws.Range("F72:V72").Value2 = $rowB
$ws.Range("F73:V73").Value2 = $rowA

$rowA = $ws.Range("F82:V82").Value2
$rowB = $ws.Range("F83:V83").Value2
$ws.Range("F82:V82").Value2 = $rowB
$ws.Range("F83:V83").Value2 = $rowA

$rowA = $ws.Range("F85:V85").Value2
$rowB = $ws.Range("F87:V87").Value2
$ws.Range("F85:V85").Value2 = $rowB
$ws.Range("F87:V87").Value2 = $rowA

$rowA = $ws.Range("F93:V93").Value2
$rowB = $ws.Range("F94:V94").Value2
$ws.Range("F93:V93").Value2 = $rowB
$ws.Range("F94:V94").Value2 = $rowA

$rowA = $ws.Range("F99:V99").Value2
$rowB = $ws.Range("F100:V100").Value2
$ws.Range("F99:V99").Value2 = $rowB
$ws.Range("F100:V100").Value2 = $rowA

$rowA = $ws.Range("F107:V107").Value2
$rowB = $ws.Range("F108:V108").Value2
$ws.Range("F107:V107").Value2 = $rowB
$ws.Range("F108:V108").Value2 = $rowA

$rowA = $ws.Range("F114:V114").Value2
$rowB = $ws.Range("F115:V115").Value2
$ws.Range("F114:V114").Value2 = $rowB
$ws.Range("F115:V115").Value2 = $rowA

$rowA = $ws.Range("F130:V130").Value2
$rowB = $ws.Range("F131:V131").Value2
$ws.Range("F130:V130").Value2 = $rowB
$ws.Range("F131:V131").Value2 = $rowA

# --- Append 11 new rows (134-144), copying format from row 133 first ---
$ws.Range("A133:V133").Copy()
$ws.Range("A134:V144").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A134").Value2 = 133
$ws.Range("B134").Value2 = "spain"
$ws.Range("C134").Value2 = "laliga2"
$ws.Range("D134").Value2 = "2023-2024"
$ws.Range("E134").Value2 = 45226.875
$ws.Range("F134").Value2 = "Eibar"
$ws.Range("G134").Value2 = 5
$ws.Range("H134").Value2 = "Valladolid"
$ws.Range("I134").Value2 = 1
$ws.Range("J134").Value2 = 1.95
$ws.Range("K134").Value2 = "22/10/2023 21:12"
$ws.Range("L134").Value2 = 2.06
$ws.Range("M134").Value2 = "27/10/2023 20:58"
$ws.Range("N134").Value2 = 3.45
$ws.Range("O134").Value2 = "22/10/2023 21:12"
$ws.Range("P134").Value2 = 3.1
$ws.Range("Q134").Value2 = "27/10/2023 20:58"
$ws.Range("R134").Value2 = 4.32
$ws.Range("S134").Value2 = "22/10/2023 21:12"
$ws.Range("T134").Value2 = 4.53
$ws.Range("U134").Value2 = "27/10/2023 20:58"
$ws.Range("V134").Value2 = "https://www.betexplorer.com/football/spain/laliga2/eibar-valladolid/06F4WZL1/"
$ws.Range("A135").Value2 = 134
$ws.Range("B135").Value2 = "spain"
$ws.Range("C135").Value2 = "laliga2"
$ws.Range("D135").Value2 = "2023-2024"
$ws.Range("E135").Value2 = 45227.58333333334
$ws.Range("F135").Value2 = "Eldense"
$ws.Range("G135").Value2 = 2
$ws.Range("H135").Value2 = "Amorebieta"
$ws.Range("I135").Value2 = 0
$ws.Range("J135").Value2 = 1.97
$ws.Range("K135").Value2 = "22/10/2023 14:13"
$ws.Range("L135").Value2 = 1.76
$ws.Range("M135").Value2 = "28/10/2023 13:56"
$ws.Range("N135").Value2 = 3.35
$ws.Range("O135").Value2 = "22/10/2023 14:13"
$ws.Range("P135").Value2 = 3.77
$ws.Range("Q135").Value2 = "28/10/2023 13:56"
$ws.Range("R135").Value2 = 4.4
$ws.Range("S135").Value2 = "22/10/2023 14:13"
$ws.Range("T135").Value2 = 5.11
$ws.Range("U135").Value2 = "28/10/2023 13:56"
$ws.Range("V135").Value2 = "https://www.betexplorer.com/football/spain/laliga2/eldense-amorebieta/SY8gypPQ/"
$ws.Range("A136").Value2 = 135
$ws.Range("B136").Value2 = "spain"
$ws.Range("C136").Value2 = "laliga2"
$ws.Range("D136").Value2 = "2023-2024"
$ws.Range("E136").Value2 = 45227.77083333334
$ws.Range("F136").Value2 = "Gijon"
$ws.Range("G136").Value2 = 2
$ws.Range("H136").Value2 = "Espanyol"
$ws.Range("I136").Value2 = 0
$ws.Range("J136").Value2 = 2.67
$ws.Range("K136").Value2 = "22/10/2023 20:15"
$ws.Range("L136").Value2 = 2.53
$ws.Range("M136").Value2 = "28/10/2023 18:27"
$ws.Range("N136").Value2 = 3.19
$ws.Range("O136").Value2 = "22/10/2023 20:15"
$ws.Range("P136").Value2 = 3.2
$ws.Range("Q136").Value2 = "28/10/2023 18:26"
$ws.Range("R136").Value2 = 2.91
$ws.Range("S136").Value2 = "22/10/2023 20:15"
$ws.Range("T136").Value2 = 3.13
$ws.Range("U136").Value2 = "28/10/2023 18:27"
$ws.Range("V136").Value2 = "https://www.betexplorer.com/football/spain/laliga2/gijon-espanyol/C6pQVNnl/"
$ws.Range("A137").Value2 = 136
$ws.Range("B137").Value2 = "spain"
$ws.Range("C137").Value2 = "laliga2"
$ws.Range("D137").Value2 = "2023-2024"
$ws.Range("E137").Value2 = 45227.77083333334
$ws.Range("F137").Value2 = "Andorra"
$ws.Range("G137").Value2 = 2
$ws.Range("H137").Value2 = "Levante"
$ws.Range("I137").Value2 = 0
$ws.Range("J137").Value2 = 2.53
$ws.Range("K137").Value2 = "22/10/2023 21:12"
$ws.Range("L137").Value2 = 2.39
$ws.Range("M137").Value2 = "28/10/2023 18:29"
$ws.Range("N137").Value2 = 3.26
$ws.Range("O137").Value2 = "22/10/2023 21:12"
$ws.Range("P137").Value2 = 3.13
$ws.Range("Q137").Value2 = "28/10/2023 18:29"
$ws.Range("R137").Value2 = 3.05
$ws.Range("S137").Value2 = "22/10/2023 21:12"
$ws.Range("T137").Value2 = 3.44
$ws.Range("U137").Value2 = "28/10/2023 18:29"
$ws.Range("V137").Value2 = "https://www.betexplorer.com/football/spain/laliga2/fc-andorra-levante/WjPJa2Ye/"
$ws.Range("A138").Value2 = 137
$ws.Range("B138").Value2 = "spain"
$ws.Range("C138").Value2 = "laliga2"
$ws.Range("D138").Value2 = "2023-2024"
$ws.Range("E138").Value2 = 45227.875
$ws.Range("F138").Value2 = "Elche"
$ws.Range("G138").Value2 = 2
$ws.Range("H138").Value2 = "Tenerife"
$ws.Range("I138").Value2 = 1
$ws.Range("J138").Value2 = 2.09
$ws.Range("K138").Value2 = "22/10/2023 16:42"
$ws.Range("L138").Value2 = 2.05
$ws.Range("M138").Value2 = "28/10/2023 20:38"
$ws.Range("N138").Value2 = 3.25
$ws.Range("O138").Value2 = "22/10/2023 16:42"
$ws.Range("P138").Value2 = 3.2
$ws.Range("Q138").Value2 = "28/10/2023 20:38"
$ws.Range("R138").Value2 = 4.07
$ws.Range("S138").Value2 = "22/10/2023 16:42"
$ws.Range("T138").Value2 = 4.37
$ws.Range("U138").Value2 = "28/10/2023 20:38"
$ws.Range("V138").Value2 = "https://www.betexplorer.com/football/spain/laliga2/elche-tenerife/vcE8Vgy8/"
$ws.Range("A139").Value2 = 138
$ws.Range("B139").Value2 = "spain"
$ws.Range("C139").Value2 = "laliga2"
$ws.Range("D139").Value2 = "2023-2024"
$ws.Range("E139").Value2 = 45228.58333333334
$ws.Range("F139").Value2 = "Mirandes"
$ws.Range("G139").Value2 = 2
$ws.Range("H139").Value2 = "FC Cartagena SAD"
$ws.Range("I139").Value2 = 1
$ws.Range("J139").Value2 = 2.22
$ws.Range("K139").Value2 = "23/10/2023 21:12"
$ws.Range("L139").Value2 = 2.26
$ws.Range("M139").Value2 = "29/10/2023 13:52"
$ws.Range("N139").Value2 = 3.35
$ws.Range("O139").Value2 = "23/10/2023 21:12"
$ws.Range("P139").Value2 = 3.24
$ws.Range("Q139").Value2 = "29/10/2023 13:52"
$ws.Range("R139").Value2 = 3.44
$ws.Range("S139").Value2 = "23/10/2023 21:12"
$ws.Range("T139").Value2 = 3.61
$ws.Range("U139").Value2 = "29/10/2023 13:52"
$ws.Range("V139").Value2 = "https://www.betexplorer.com/football/spain/laliga2/mirandes-fc-cartagena-sad/YFoMW3Xs/"
$ws.Range("A140").Value2 = 139
$ws.Range("B140").Value2 = "spain"
$ws.Range("C140").Value2 = "laliga2"
$ws.Range("D140").Value2 = "2023-2024"
$ws.Range("E140").Value2 = 45228.67708333334
$ws.Range("F140").Value2 = "Huesca"
$ws.Range("G140").Value2 = 0
$ws.Range("H140").Value2 = "Albacete"
$ws.Range("I140").Value2 = 0
$ws.Range("J140").Value2 = 2.93
$ws.Range("K140").Value2 = "22/10/2023 20:15"
$ws.Range("L140").Value2 = 3.3
$ws.Range("M140").Value2 = "29/10/2023 16:06"
$ws.Range("N140").Value2 = 2.94
$ws.Range("O140").Value2 = "22/10/2023 20:15"
$ws.Range("P140").Value2 = 2.79
$ws.Range("Q140").Value2 = "29/10/2023 16:06"
$ws.Range("R140").Value2 = 2.8
$ws.Range("S140").Value2 = "22/10/2023 20:15"
$ws.Range("T140").Value2 = 2.73
$ws.Range("U140").Value2 = "29/10/2023 16:06"
$ws.Range("V140").Value2 = "https://www.betexplorer.com/football/spain/laliga2/huesca-albacete/K85ow6fE/"
$ws.Range("A141").Value2 = 140
$ws.Range("B141").Value2 = "spain"
$ws.Range("C141").Value2 = "laliga2"
$ws.Range("D141").Value2 = "2023-2024"
$ws.Range("E141").Value2 = 45228.77083333334
$ws.Range("F141").Value2 = "Leganes"
$ws.Range("G141").Value2 = 1
$ws.Range("H141").Value2 = "Villarreal B"
$ws.Range("I141").Value2 = 0
$ws.Range("J141").Value2 = 1.69
$ws.Range("K141").Value2 = "22/10/2023 20:15"
$ws.Range("L141").Value2 = 1.8
$ws.Range("M141").Value2 = "29/10/2023 18:29"
$ws.Range("N141").Value2 = 3.91
$ws.Range("O141").Value2 = "22/10/2023 20:15"
$ws.Range("P141").Value2 = 3.67
$ws.Range("Q141").Value2 = "29/10/2023 18:29"
$ws.Range("R141").Value2 = 5.32
$ws.Range("S141").Value2 = "22/10/2023 20:15"
$ws.Range("T141").Value2 = 4.96
$ws.Range("U141").Value2 = "29/10/2023 18:29"
$ws.Range("V141").Value2 = "https://www.betexplorer.com/football/spain/laliga2/leganes-villarreal/tI6svnu8/"
$ws.Range("A142").Value2 = 141
$ws.Range("B142").Value2 = "spain"
$ws.Range("C142").Value2 = "laliga2"
$ws.Range("D142").Value2 = "2023-2024"
$ws.Range("E142").Value2 = 45228.77083333334
$ws.Range("F142").Value2 = "R. Oviedo"
$ws.Range("G142").Value2 = 2
$ws.Range("H142").Value2 = "Alcorcon"
$ws.Range("I142").Value2 = 0
$ws.Range("J142").Value2 = 2.02
$ws.Range("K142").Value2 = "23/10/2023 21:12"
$ws.Range("L142").Value2 = 1.83
$ws.Range("M142").Value2 = "29/10/2023 18:01"
$ws.Range("N142").Value2 = 3.14
$ws.Range("O142").Value2 = "23/10/2023 21:12"
$ws.Range("P142").Value2 = 3.3
$ws.Range("Q142").Value2 = "29/10/2023 18:27"
$ws.Range("R142").Value2 = 4.58
$ws.Range("S142").Value2 = "23/10/2023 21:12"
$ws.Range("T142").Value2 = 5.59
$ws.Range("U142").Value2 = "29/10/2023 18:27"
$ws.Range("V142").Value2 = "https://www.betexplorer.com/football/spain/laliga2/r-oviedo-alcorcon/EV4kxQ9K/"
$ws.Range("A143").Value2 = 142
$ws.Range("B143").Value2 = "spain"
$ws.Range("C143").Value2 = "laliga2"
$ws.Range("D143").Value2 = "2023-2024"
$ws.Range("E143").Value2 = 45228.875
$ws.Range("F143").Value2 = "Racing Santander"
$ws.Range("G143").Value2 = 1
$ws.Range("H143").Value2 = "Ferrol"
$ws.Range("I143").Value2 = 3
$ws.Range("J143").Value2 = 2.22
$ws.Range("K143").Value2 = "22/10/2023 22:13"
$ws.Range("L143").Value2 = 2.18
$ws.Range("M143").Value2 = "29/10/2023 20:33"
$ws.Range("N143").Value2 = 3.24
$ws.Range("O143").Value2 = "22/10/2023 22:13"
$ws.Range("P143").Value2 = 3.19
$ws.Range("Q143").Value2 = "29/10/2023 20:49"
$ws.Range("R143").Value2 = 3.68
$ws.Range("S143").Value2 = "22/10/2023 22:13"
$ws.Range("T143").Value2 = 3.92
$ws.Range("U143").Value2 = "29/10/2023 20:56"
$ws.Range("V143").Value2 = "https://www.betexplorer.com/football/spain/laliga2/racing-santander-ferrol/ttkUUs2f/"
$ws.Range("A144").Value2 = 143
$ws.Range("B144").Value2 = "spain"
$ws.Range("C144").Value2 = "laliga2"
$ws.Range("D144").Value2 = "2023-2024"
$ws.Range("E144").Value2 = 45229.875
$ws.Range("F144").Value2 = "Burgos CF"
$ws.Range("G144").Value2 = 1
$ws.Range("H144").Value2 = "Zaragoza"
$ws.Range("I144").Value2 = 1
$ws.Range("J144").Value2 = 2.4
$ws.Range("K144").Value2 = "23/10/2023 22:12"
$ws.Range("L144").Value2 = 2.44
$ws.Range("M144").Value2 = "30/10/2023 20:56"
$ws.Range("N144").Value2 = 3.09
$ws.Range("O144").Value2 = "23/10/2023 22:12"
$ws.Range("P144").Value2 = 2.84
$ws.Range("Q144").Value2 = "30/10/2023 20:55"
$ws.Range("R144").Value2 = 3.43
$ws.Range("S144").Value2 = "23/10/2023 22:12"
$ws.Range("T144").Value2 = 3.76
$ws.Range("U144").Value2 = "30/10/2023 20:56"
$ws.Range("V144").Value2 = "https://www.betexplorer.com/football/spain/laliga2/burgos-cf-zaragoza/6qG0XF6e/"

